# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 18:22"

# --- Swap Cuba / Nueva Zelanda ranking (Cuba overtakes Nueva Zelanda) ---
# Row 80 keeps Nueva Zelanda's old figures but the country becomes Cuba with
# its newly updated figures; row 81 becomes Nueva Zelanda with the figures
# that used to sit in row 80.
$ws.Range("A80").Value = "Cuba"
$ws.Range("B80").Value = 1501
$ws.Range("C80").Value = 34
$ws.Range("D80").Value = 681
$ws.Range("E80").Value = 759
$ws.Range("F80").Value = 10
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = 61

$ws.Range("A81").Value = "Nueva Zelanda"
$ws.Range("B81").Value = 1476
$ws.Range("C81").Value = 2
$ws.Range("D81").Value = 1241
$ws.Range("E81").Value = 216
$ws.Range("F81").Value = 1
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 19

# --- Swap Somalia / Guatemala ranking (Somalia overtakes Guatemala) ---
$ws.Range("A106").Value = "Somalia"
$ws.Range("B106").Value = 601
$ws.Range("C106").Value = 19
$ws.Range("D106").Value = 31
$ws.Range("E106").Value = 542
$ws.Range("F106").Value = 2
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 28

$ws.Range("A107").Value = "Guatemala"
$ws.Range("B107").Value = 585
$ws.Range("C107").Value = 28
$ws.Range("D107").Value = 65
$ws.Range("E107").Value = 504
$ws.Range("F107").Value = 5
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 16

# --- Plain numeric refreshes (no reordering) ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 1069294
$ws.Range("C4").Value = 5100
$ws.Range("D4").Value = 148291
$ws.Range("E4").Value = 858795
$ws.Range("G4").Value = 553
$ws.Range("H4").Value = 62208

# Italia (row 6)
$ws.Range("B6").Value = 205463
$ws.Range("C6").Value = 1872
$ws.Range("D6").Value = 75945
$ws.Range("E6").Value = 101551
$ws.Range("F6").Value = 1694
$ws.Range("G6").Value = 285
$ws.Range("H6").Value = 27967

# Turquia (row 10)
$ws.Range("B10").Value = 120204
$ws.Range("C10").Value = 2615
$ws.Range("D10").Value = 48886
$ws.Range("E10").Value = 68144
$ws.Range("F10").Value = 1514
$ws.Range("G10").Value = 93
$ws.Range("H10").Value = 3174

# Singapur (row 27)
$ws.Range("D27").Value = 1244
$ws.Range("E27").Value = 14910
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 15

# Islandia (row 75)
$ws.Range("D75").Value = 1670
$ws.Range("E75").Value = 117

# Sri Lanka (row 103)
$ws.Range("B103").Value = 660
$ws.Range("C103").Value = 11
$ws.Range("E103").Value = 514
